$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new booking rows for room_id 6 and 7 (check-in allowing many rooms)
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 6
$ws.Range("C15").Value = 44905

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 7
$ws.Range("C16").Value = 44905

# Match the date formatting used by the existing booking_date column
$ws.Range("C14").Copy()
$ws.Range("C15:C16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Move the active selection below the newly added rows
$ws.Range("C17").Select()
